$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '31.154.73'
$ws.Range("E2").Value = '  +1.82%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '1.988.54'
$ws.Range("E3").Value = '  +5.64%  '

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.06%  '

# Row 5: XRP
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7925'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +67.56%  '

# Row 6: BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '254.32'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.21%  '

# Row 7: USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9998'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.00%  '

# Row 8: Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3506'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +21.37%  '

# Row 9: Solana
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '28.03'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +26.05%  '

# Row 10: Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06993'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +6.99%  '

# Row 11: Polygon
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8436'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +8.98%  '

# Row 12: TRON
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08186'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +4.59%  '

# Row 13: WrappedEther
$ws.Range("D13").Value = '1.988.51'
$ws.Range("E13").Value = '  +5.67%  '

# Row 14: Litecoin
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '100.23'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.63%  '

# Row 15: Polkadot
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.582'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +6.30%  '

# Row 16: Avalanche
$ws.Range("E16").Value = '  +16.13%  '

# Row 17: BitcoinCash
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '272.97'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.54%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = '31.156.83'
$ws.Range("E18").Value = '  +1.91%  '

# Row 19: Uniswap
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.861'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +9.20%  '

# Row 20: ShibaInu
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007904'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +5.03%  '

# Row 21: WrappedliquidstakedEther2.0
$ws.Range("D21").Value = '2.251.34'
$ws.Range("E21").Value = '  +5.98%  '

# Row 22: Dai
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9996'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.01%  '

# Row 23: BinanceUSD
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9984'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.16%  '

# Row 24: Chainlink
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.044'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +9.90%  '

# Row 25: Cosmos
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.00'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +9.52%  '

# Row 26: Stellar
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1508'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +55.34%  '

# Row 27: Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '165.19'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.29%  '

# Row 28: EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.95'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +4.42%  '

# Row 29: LidoDAOToken
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.314'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +20.86%  '

# Row 30: PancakeSwap
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.594'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +5.96%  '

# Row 31: Toncoin
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.357'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.86%  '

# Row 32: Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.580'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +7.38%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.403'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +4.98%  '

# Row 34: Hedera
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05221'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +7.64%  '

# Row 35: ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.226'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +8.45%  '

# Row 36: ImmutableX
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7779'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +11.58%  '

# Row 37: HuobiToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.760'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.55%  '

# Row 38: Frax
$ws.Range("B38").Value = 'Frax'
$ws.Range("C38").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9989'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.07%  '

# Row 39: VeChain
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02002'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +4.35%  '

# Row 40: MXToken
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.890'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.35%  '

# Row 41: FraxShare
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.617'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +5.22%  '

# Row 42: Aave
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '79.12'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.80%  '

# Row 43: TheSandbox
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4658'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +9.49%  '

# Row 44: RenderToken
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.121'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +6.62%  '

# Row 45: Quant
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '105.05'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +3.58%  '

# Row 46: TrustWalletToken
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.8476'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.52%  '

# Row 47: PaxDollar
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9996'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.01%  '

# Row 48: Aptos
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.665'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +8.84%  '

# Row 49: EnergySwap
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.854'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.09%  '

# Row 50: Elrond
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.73'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +4.28%  '

# Row 51: Decentraland
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4295'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +8.65%  '
